# Insert one new data row right before the current row 327 (the "Paine - 1a (guarda)"
# entry dated 44515). This pushes that row and every row below it down by one
# position (327 -> 328, 328 -> 329, ... 435 -> 436) and enlarges the used range
# from A1:R435 to A1:R436.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(327).Insert()

# Populate the freshly inserted row 327 with the new record.
$ws.Cells.Item(327, 1).Value  = 4
$ws.Cells.Item(327, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(327, 3).Value  = "Los Lagos"
$ws.Cells.Item(327, 4).Value  = 44988
$ws.Cells.Item(327, 5).Value  = 10
$ws.Cells.Item(327, 6).Value  = 100112045
$ws.Cells.Item(327, 7).Value  = "Zapallo"
$ws.Cells.Item(327, 8).Value  = "Paine"
$ws.Cells.Item(327, 9).Value  = "1a (cosecha)"
$ws.Cells.Item(327, 10).Value = 1200
$ws.Cells.Item(327, 11).Value = 500
$ws.Cells.Item(327, 12).Value = 600
$ws.Cells.Item(327, 13).Value = 550
$ws.Cells.Item(327, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(327, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(327, 16).Value = 550
$ws.Cells.Item(327, 17).Value = 1
$ws.Cells.Item(327, 18).Value = "Hortaliza"
